# #5: property aircraft done
#
# The "property_category" column on the 建物 (building) and 汽車 (car)
# sheets was left as "land" (copy/paste leftover from the 土地 sheet).
# Correct it to the proper category for each sheet.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: column I is property_category, data rows 2-4
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I4").Value2 = "building"

# 汽車 (car) sheet: column H is property_category, data rows 2-4
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2:H4").Value2 = "car"
